$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 552, shifting existing rows 552:664 down to 553:665.
$ws.Rows.Item(552).Insert()

# Populate the newly inserted row 552 with its data.
$ws.Range("A552").Value = 3
$ws.Range("B552").Value = "Femacal de La Calera"
$ws.Range("C552").Value = "Coquimbo"
$ws.Range("D552").Value = 45258
$ws.Range("E552").Value = 5
$ws.Range("F552").Value = 100114013
$ws.Range("G552").Value = "Zanahoria"
$ws.Range("H552").Value = "Sin especificar"
$ws.Range("I552").Value = "Primera"
$ws.Range("J552").Value = 160
$ws.Range("K552").Value = 6000
$ws.Range("L552").Value = 6000
$ws.Range("M552").Value = 6000
$ws.Range("N552").Value = "$/saco 20 kilos"
$ws.Range("O552").Value = "Provincia de Quillota"
$ws.Range("P552").Value = 300
$ws.Range("Q552").Value = 20
$ws.Range("R552").Value = "Hortaliza"
